$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.412870246130808
$ws.Range("C2").Value = 0.1682565082276142
$ws.Range("D2").Value = 0.09808484263608008
$ws.Range("F2").Value = 2.069221938363427
$ws.Range("G2").Value = 1.42772883144643
$ws.Range("H2").Value = 1.317309718956423
$ws.Range("I2").Value = 1.252771063368982
$ws.Range("J2").Value = 0.1774909734679824
$ws.Range("L2").Value = 0.3717099617930302
$ws.Range("M2").Value = 0.356788751973113

# Row 3
$ws.Range("B3").Value = 1.317231774682284
$ws.Range("C3").Value = 0.1484401673830291
$ws.Range("D3").Value = 0.09788900112230436
$ws.Range("F3").Value = 2.076900983398119
$ws.Range("G3").Value = 1.430531225437946
$ws.Range("H3").Value = 1.324945147423506
$ws.Range("I3").Value = 1.264583940846315
$ws.Range("J3").Value = 0.178386177027777
$ws.Range("L3").Value = 0.3688572064345621
$ws.Range("M3").Value = 0.3416913777316068

# Row 4
$ws.Range("B4").Value = 1.258864914646665
$ws.Range("C4").Value = 0.1362136664180582
$ws.Range("D4").Value = 0.09777671807255039
$ws.Range("F4").Value = 2.08282140055762
$ws.Range("G4").Value = 1.433253289970764
$ws.Range("H4").Value = 1.330322260059575
$ws.Range("I4").Value = 1.272562690787449
$ws.Range("J4").Value = 0.1789678888291739
$ws.Range("L4").Value = 0.3672382464765747
$ws.Range("M4").Value = 0.3325395810519183

# Row 5
$ws.Range("B5").Value = 1.235170557281094
$ws.Range("C5").Value = 0.131216468386441
$ws.Range("D5").Value = 0.09773297431126693
$ws.Range("F5").Value = 2.085536886252015
$ws.Range("G5").Value = 1.434613879506969
$ws.Range("H5").Value = 1.332686629165636
$ws.Range("I5").Value = 1.275996357143416
$ws.Range("J5").Value = 0.1792130208201241
$ws.Range("L5").Value = 0.3666119576425899
$ws.Range("M5").Value = 0.3288400518330832

# Row 6
$ws.Range("B6").Value = 1.23124163591774
$ws.Range("C6").Value = 0.1303857961307813
$ws.Range("D6").Value = 0.09772583253457867
$ws.Range("F6").Value = 2.086006075663171
$ws.Range("G6").Value = 1.434854969808129
$ws.Range("H6").Value = 1.333089686356075
$ws.Range("I6").Value = 1.27657751720967
$ws.Range("J6").Value = 0.1792542133751196
$ws.Range("L6").Value = 0.3665099860360925
$ws.Range("M6").Value = 0.3282275602957938

# Row 7
$ws.Range("B7").Value = 1.258544995761213
$ws.Range("C7").Value = 0.1361463322675718
$ws.Range("D7").Value = 0.09777611996615931
$ws.Range("F7").Value = 2.082856796590313
$ws.Range("G7").Value = 1.433270622388406
$ws.Range("H7").Value = 1.330353445800299
$ws.Range("I7").Value = 1.27260826070998
$ws.Range("J7").Value = 0.1789711620244638
$ws.Range("L7").Value = 0.3672296645576907
$ws.Range("M7").Value = 0.332489566544389

# Row 8
$ws.Range("B8").Value = 1.379820994680188
$ws.Range("C8").Value = 0.1614361619494957
$ws.Range("D8").Value = 0.09801566915283999
$ws.Range("F8").Value = 2.071619325341942
$ws.Range("G8").Value = 1.428486967548267
$ws.Range("H8").Value = 1.319799375187998
$ws.Range("I8").Value = 1.256693428338153
$ws.Range("J8").Value = 0.1777929988566065
$ws.Range("L8").Value = 0.3706988574959666
$ws.Range("M8").Value = 0.3515587997321745

# Row 9
$ws.Range("B9").Value = 1.620424345781544
$ws.Range("C9").Value = 0.2105585844390134
$ws.Range("D9").Value = 0.09854824272261453
$ws.Range("F9").Value = 2.059160684754644
$ws.Range("G9").Value = 1.427075150102027
$ws.Range("H9").Value = 1.304573764590742
$ws.Range("I9").Value = 1.231250782770047
$ws.Range("J9").Value = 0.1757360589786892
$ws.Range("L9").Value = 0.3785508749060824
$ws.Range("M9").Value = 0.3898829260926036

# Row 10
$ws.Range("B10").Value = 1.798857782007588
$ws.Range("C10").Value = 0.2463637394680518
$ws.Range("D10").Value = 0.09897739402687833
$ws.Range("F10").Value = 2.055867417785137
$ws.Range("G10").Value = 1.430930851593473
$ws.Range("H10").Value = 1.296730072474048
$ws.Range("I10").Value = 1.216086113619951
$ws.Range("J10").Value = 0.1743780816938028
$ws.Range("L10").Value = 0.384955515285526
$ws.Range("M10").Value = 0.4185998269213371

# Row 11
$ws.Range("B11").Value = 1.880387075911585
$ws.Range("C11").Value = 0.2625912050398256
$ws.Range("D11").Value = 0.09918076761935879
$ws.Range("F11").Value = 2.055646438715343
$ws.Range("G11").Value = 1.433755054189149
$ws.Range("H11").Value = 1.293889267515652
$ws.Range("I11").Value = 1.209956185627753
$ws.Range("J11").Value = 0.17379331812743
$ws.Range("L11").Value = 0.3880064994858401
$ws.Range("M11").Value = 0.4317844118928704

# Row 12
$ws.Range("B12").Value = 1.911310853353086
$ws.Range("C12").Value = 0.2687274055700186
$ws.Range("D12").Value = 0.09925894419768255
$ws.Range("F12").Value = 2.055746756628508
$ws.Range("G12").Value = 1.434978981998825
$ws.Range("H12").Value = 1.292918244972753
$ws.Range("I12").Value = 1.207745691851365
$ws.Range("J12").Value = 0.1735766072120875
$ws.Range("L12").Value = 0.389181519763909
$ws.Range("M12").Value = 0.4367943282525175

# Row 13
$ws.Range("B13").Value = 1.904648642198026
$ws.Range("C13").Value = 0.2674062572997968
$ws.Range("D13").Value = 0.09924205586441204
$ws.Range("F13").Value = 2.055716962491474
$ws.Range("G13").Value = 1.43470850840967
$ws.Range("H13").Value = 1.293122711981709
$ws.Range("I13").Value = 1.208216830263488
$ws.Range("J13").Value = 0.1736230698501249
$ws.Range("L13").Value = 0.3889275843407205
$ws.Range("M13").Value = 0.4357145915694787

# Row 14
$ws.Range("B14").Value = 1.882930190489162
$ws.Range("C14").Value = 0.2630962106003665
$ws.Range("D14").Value = 0.09918717598493743
$ws.Range("F14").Value = 2.055651002352846
$ws.Range("G14").Value = 1.433852648540025
$ws.Range("H14").Value = 1.293807281186048
$ws.Range("I14").Value = 1.209772105351576
$ws.Range("J14").Value = 0.1737753945643945
$ws.Range("L14").Value = 0.3881027752362058
$ws.Range("M14").Value = 0.432196237243204

# Row 15
$ws.Range("B15").Value = 1.869633539654046
$ws.Range("C15").Value = 0.2604550353494517
$ws.Range("D15").Value = 0.09915371173382326
$ws.Range("F15").Value = 2.05563457210846
$ws.Range("G15").Value = 1.433348542196214
$ws.Range("H15").Value = 1.294240242509233
$ws.Range("I15").Value = 1.210739190387088
$ws.Range("J15").Value = 0.1738693128559703
$ws.Range("L15").Value = 0.3876001159223392
$ws.Range("M15").Value = 0.4300433772039653

# Row 16
$ws.Range("B16").Value = 1.793536770402682
$ws.Range("C16").Value = 0.2453020115416109
$ws.Range("D16").Value = 0.09896426635085831
$ws.Range("F16").Value = 2.055907588258179
$ws.Range("G16").Value = 1.430767867392532
$ws.Range("H16").Value = 1.296930373394389
$ws.Range("I16").Value = 1.21650220631053
$ws.Range("J16").Value = 0.17441695968202
$ws.Range("L16").Value = 0.3847588849768471
$ws.Range("M16").Value = 0.4177406019580445

# Row 17
$ws.Range("B17").Value = 1.746944946221504
$ws.Range("C17").Value = 0.2359905660701429
$ws.Range("D17").Value = 0.09885012907035318
$ws.Range("F17").Value = 2.056402428093321
$ws.Range("G17").Value = 1.429459220231777
$ws.Range("H17").Value = 1.298767067104578
$ws.Range("I17").Value = 1.220234660003669
$ws.Range("J17").Value = 0.1747613597920346
$ws.Range("L17").Value = 0.3830510290931954
$ws.Range("M17").Value = 0.4102241148562058

# Row 18
$ws.Range("B18").Value = 1.720180426844706
$ws.Range("C18").Value = 0.2306291816673536
$ws.Range("D18").Value = 0.09878524799164623
$ws.Range("F18").Value = 2.056807236997301
$ws.Range("G18").Value = 1.428807221702428
$ws.Range("H18").Value = 1.299891936055175
$ws.Range("I18").Value = 1.22245379321344
$ws.Range("J18").Value = 0.174962555540521
$ws.Range("L18").Value = 0.3820816613085753
$ws.Range("M18").Value = 0.4059122437203655

# Row 19
$ws.Range("B19").Value = 1.711124274511235
$ws.Range("C19").Value = 0.2288129350737904
$ws.Range("D19").Value = 0.09876341248709863
$ws.Range("F19").Value = 2.056964930588421
$ws.Range("G19").Value = 1.428603744395645
$ws.Range("H19").Value = 1.300284549707399
$ws.Range("I19").Value = 1.223217567512926
$ws.Range("J19").Value = 0.175031210934355
$ws.Range("L19").Value = 0.381755676261804
$ws.Range("M19").Value = 0.4044542851753832

# Row 20
$ws.Range("B20").Value = 1.751901229410407
$ws.Range("C20").Value = 0.2369823746191457
$ws.Range("D20").Value = 0.0988621997919239
$ws.Range("F20").Value = 2.056337309811113
$ws.Range("G20").Value = 1.429588101702791
$ws.Range("H20").Value = 1.298564462457222
$ws.Range("I20").Value = 1.219829846731848
$ws.Range("J20").Value = 0.174724376486795
$ws.Range("L20").Value = 0.3832314940926693
$ws.Range("M20").Value = 0.4110230780738462

# Row 21
$ws.Range("B21").Value = 1.889308067962304
$ws.Range("C21").Value = 0.2643624150292965
$ws.Range("D21").Value = 0.09920326402889756
$ws.Range("F21").Value = 2.055665380030959
$ws.Range("G21").Value = 1.434099838715255
$ws.Range("H21").Value = 1.29360336319408
$ws.Range("I21").Value = 1.209312274832584
$ws.Range("J21").Value = 0.1737305249670182
$ws.Range("L21").Value = 0.3883445083066732
$ws.Range("M21").Value = 0.4332291981031986

# Row 22
$ws.Range("B22").Value = 1.979404612979465
$ws.Range("C22").Value = 0.2822055873604654
$ws.Range("D22").Value = 0.09943294666161862
$ws.Range("F22").Value = 2.05629887701231
$ws.Range("G22").Value = 1.437949114865233
$ws.Range("H22").Value = 1.290971502255701
$ws.Range("I22").Value = 1.203084247155125
$ws.Range("J22").Value = 0.1731085265586678
$ws.Range("L22").Value = 0.3918008201733301
$ws.Range("M22").Value = 0.4478423383117089

# Row 23
$ws.Range("B23").Value = 1.931291960576118
$ws.Range("C23").Value = 0.2726870710365858
$ws.Range("D23").Value = 0.09930974336367271
$ws.Range("F23").Value = 2.055862506319158
$ws.Range("G23").Value = 1.435812091915807
$ws.Range("H23").Value = 1.292320268417768
$ws.Range("I23").Value = 1.206349084850167
$ws.Range("J23").Value = 0.1734379844224119
$ws.Range("L23").Value = 0.3899456597087863
$ws.Range("M23").Value = 0.440033936970579

# Row 24
$ws.Range("B24").Value = 1.749660425967932
$ws.Range("C24").Value = 0.2365340032345102
$ws.Range("D24").Value = 0.09885674032025449
$ws.Range("F24").Value = 2.056366375015784
$ws.Range("G24").Value = 1.429529521809982
$ws.Range("H24").Value = 1.298655845308218
$ws.Range("I24").Value = 1.220012634510525
$ws.Range("J24").Value = 0.1747410866836532
$ws.Range("L24").Value = 0.3831498669199505
$ws.Range("M24").Value = 0.4106618373282203

# Row 25
$ws.Range("B25").Value = 1.555040519400279
$ws.Range("C25").Value = 0.197319883378384
$ws.Range("D25").Value = 0.09839748733345743
$ws.Range("F25").Value = 2.061503509873077
$ws.Range("G25").Value = 1.426600350576052
$ws.Range("H25").Value = 1.308106242302273
$ws.Range("I25").Value = 1.237515055009183
$ws.Range("J25").Value = 0.1762655132143429
$ws.Range("L25").Value = 0.3763147859855991
$ws.Range("M25").Value = 0.3794164101658382

Write-Output "applied 380 kV case values"